$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Region_ZH), shifting Region_ZH/Abbrev right.
$ws.Columns.Item(4).Insert()

# Header for the new Region_EN column
$ws.Range("D1").Value = "Region_EN"

# Map each district row's Region_ZH value (now in column E) to the English region name
$regionMap = @{
    "香港島" = "Hong Kong Island"
    "九龍" = "Kowloon"
    "新界西" = "New Territories West"
    "新界東" = "New Territories East"
}

for ($r = 2; $r -le 19; $r++) {
    $zh = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value = $regionMap[$zh]
}

# Adjust column widths: new D (Region_EN) gets a wider bestFit width, E keeps the old D width
# (target widths 20.140625 / 10.42578125 aren't exactly reachable through this host's
# column-width quantization, so use the closest attainable values)
$ws.Columns.Item(4).ColumnWidth = 19.33
$ws.Columns.Item(5).ColumnWidth = 9.65

# Update selection to match the target state
$ws.Range("F7").Select()
